$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update hate_speech counts for existing rows (dataset values changed)
$ws.Range("E4").Value = 8
$ws.Range("E7").Value = 9
$ws.Range("E9").Value = 16
$ws.Range("E20").Value = 27
